# Fix docs, typos, excel parsing and parsing test
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Typo fix: "Дата предоставление кредита" -> "Дата предоставления кредита"
$ws.Range("E1").Value = "Дата предоставления кредита"

# The sheet had a spurious, unused trailing column (F) and a bunch of
# trailing blank rows (3-11) that only existed because of formatting leaking
# into otherwise-empty cells. Drop the empty column and blank out the
# trailing rows completely (no formatting left behind).
$ws.Columns.Item(6).Delete()
$ws.Range("A3:E11").Clear()

# The populated header/data rows (1-2) were accidentally styled with white
# font text (invisible on a white background) - restore normal/automatic
# text color while keeping the same font (Calibri 11).
$ws.Range("A1:E2").Font.Color = 0

$ws.Range("E1").Select()
